$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the bottom of the data table: 2020-06-05 .. 2020-06-07
$newRows = @(
    @{ Row = 80; Date = 43987; Vals = @(57,86,132,304,161,269,123,225,181,258,184,202,159,189,115,172,53,148,3,9) },
    @{ Row = 81; Date = 43988; Vals = @(58,87,136,310,162,273,124,227,183,262,185,204,159,190,115,172,53,150,3,9) },
    @{ Row = 82; Date = 43989; Vals = @(60,88,137,314,164,275,124,231,186,263,185,206,159,191,115,172,53,151,3,9) }
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Date

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$rowNum").Value = $r.Vals[$i]
    }

    $ws.Range("V$rowNum").Formula = "=SUM(B$rowNum,D$rowNum,F$rowNum,H$rowNum,J$rowNum,L$rowNum,N$rowNum,P$rowNum,R$rowNum,T$rowNum)"
    $ws.Range("W$rowNum").Formula = "=SUM(C$rowNum,E$rowNum,G$rowNum,I$rowNum,K$rowNum,M$rowNum,O$rowNum,Q$rowNum,S$rowNum,U$rowNum)"
    $ws.Range("X$rowNum").Formula = "=(V$rowNum/(V$rowNum+W$rowNum))*100"
    $ws.Range("Y$rowNum").Formula = "=(W$rowNum/(V$rowNum+W$rowNum))*100"
}

# Update the frozen-pane view / selection to match where the user ended up
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 59
$ws.Range("V81:Y82").Select() | Out-Null
